$p = $ppt.ActivePresentation
$s = $p.Slides.Item(18)
$shp = $s.Shapes.Item(2)
$shp.TextFrame.AutoSize = 2
$tr = $shp.TextFrame.TextRange
$tr.Font.Size = 18
